$p = $ppt.ActivePresentation

# Remove the redundant "x = 1 / print(x)" example slide (SlideID 256).
# The discussion of namespaces/binding continues on the slide that used
# to follow it, so this intro slide is no longer needed.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 256) {
        $s.Delete()
        break
    }
}
